$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.108.96'
$ws.Range('E2').Value = '  +0.29%  '

$ws.Range('D3').Value = '2.322.65'
$ws.Range('E3').Value = '  +2.64%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').Value = '''253.86'
$ws.Range('E5').Value = '  +0.31%  '

$ws.Range('E6').Value = '  +1.78%  '

$ws.Range('D7').Value = '''76.47'
$ws.Range('E7').Value = '  +7.17%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').Value = '''0.655'
$ws.Range('E9').Value = '  -3.24%  '

$ws.Range('D10').Value = '''40.14'
$ws.Range('E10').Value = '  +1.17%  '

$ws.Range('D11').Value = '''0.0988'
$ws.Range('E11').Value = '  +0.92%  '

$ws.Range('D12').Value = '''7.58'
$ws.Range('E12').Value = '  -0.44%  '

$ws.Range('D13').Value = '''0.107'
$ws.Range('E13').Value = '  +2.21%  '

$ws.Range('D14').Value = '2.671.96'
$ws.Range('E14').Value = '  +2.68%  '

$ws.Range('D15').Value = '''15.48'
$ws.Range('E15').Value = '  +4.04%  '

$ws.Range('D16').Value = '''0.883'
$ws.Range('E16').Value = '  -0.82%  '

$ws.Range('D17').Value = '2.331.58'
$ws.Range('E17').Value = '  +2.57%  '

$ws.Range('D18').Value = '43.082.35'
$ws.Range('E18').Value = '  +0.33%  '

$ws.Range('E19').Value = '  +2.46%  '

$ws.Range('E20').Value = '  +0.60%  '

$ws.Range('D21').Value = '''73.10'
$ws.Range('E21').Value = '  -0.30%  '

$ws.Range('D22').Value = '''238.51'
$ws.Range('E22').Value = '  +0.24%  '

$ws.Range('D23').Value = '''2.23'
$ws.Range('E23').Value = '  +4.66%  '

$ws.Range('D24').Value = '''3.92'
$ws.Range('E24').Value = '  -1.04%  '

$ws.Range('D25').Value = '''11.64'
$ws.Range('E25').Value = '  -1.34%  '

$ws.Range('E26').Value = '  +0.09%  '

$ws.Range('D27').Value = '''2.44'
$ws.Range('E27').Value = '  -0.75%  '

$ws.Range('E28').Value = '  +1.73%  '

$ws.Range('E29').Value = '  +0.17%  '

$ws.Range('D30').Value = '''167.54'
$ws.Range('E30').Value = '  -0.26%  '

$ws.Range('D31').Value = '''0.0850'
$ws.Range('E31').Value = '  +10.07%  '

$ws.Range('D32').Value = '''6.31'
$ws.Range('E32').Value = '  -0.08%  '

$ws.Range('E33').Value = '  +1.01%  '

$ws.Range('D34').Value = '''30.53'
$ws.Range('E34').Value = '  +5.28%  '

$ws.Range('E35').Value = '  +1.90%  '

$ws.Range('E36').Value = '  +10.63%  '

$ws.Range('E37').Value = '  +2.71%  '

$ws.Range('D38').Value = '''0.0316'
$ws.Range('E38').Value = '  -2.22%  '

$ws.Range('D39').Value = '''14.06'
$ws.Range('E39').Value = '  +15.05%  '

$ws.Range('E40').Value = '  +1.93%  '

$ws.Range('D41').Value = '''5.93'
$ws.Range('E41').Value = '  +1.27%  '

$ws.Range('E42').Value = '  +8.87%  '

$ws.Range('B43').Value = 'MultiversX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D43').Value = '''62.81'
$ws.Range('E43').Value = '  -2.43%  '

$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '''9.23'
$ws.Range('E44').Value = '  +3.17%  '

$ws.Range('D45').Value = '''4.92'
$ws.Range('E45').Value = '  -2.48%  '

$ws.Range('D46').Value = '''106.27'
$ws.Range('E46').Value = '  +12.13%  '

$ws.Range('E47').Value = '  -0.41%  '

$ws.Range('E48').Value = '  -0.33%  '

$ws.Range('E49').Value = '  -0.09%  '

$ws.Range('E50').Value = '  -0.68%  '

$ws.Range('E51').Value = '  -0.66%  '
